# Remove the leading "<number><tab>" run pair that prefixes each
# Heading 1 paragraph (the "SectionNumber"-styled run plus the tab
# run), leaving only the heading's plain-text run behind.
#
#   <w:r><w:rPr><w:rStyle w:val="SectionNumber"/></w:rPr><w:t>1</w:t></w:r>
#   <w:r><w:tab/></w:r>
#   <w:r><w:t>Цель работы</w:t></w:r>
#         ->
#   <w:r><w:t>Цель работы</w:t></w:r>

$d = $word.ActiveDocument

# Collect the Heading 1 paragraphs first (collection is live, so do not
# mutate the document while iterating it).
$headingParas = @()
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Heading 1") {
        $headingParas += $p
    }
}

# Walk the headings back-to-front so deleting text in one heading does
# not invalidate the character offsets already captured for the ones
# that follow it in the document.
for ($i = $headingParas.Count - 1; $i -ge 0; $i--) {
    $p = $headingParas[$i]
    $rng = $p.Range
    $txt = $rng.Text
    $tabIndex = $txt.IndexOf([char]9)
    if ($tabIndex -ge 0) {
        $delRange = $d.Range($rng.Start, $rng.Start + $tabIndex + 1)
        $delRange.Delete()
    }
}
